$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $origStyle = $rng.Style
    $rng.Value = "'" + $val
    $rng.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "25.996.00"
$ws.Range("E2").Value = "  -0.75%  "

Set-TextValue $ws.Range("D3") "1.640.38"
$ws.Range("E3").Value = "  -1.06%  "

Set-TextValue $ws.Range("D4") "1.008"
$ws.Range("E4").Value = "  +0.07%  "

Set-TextValue $ws.Range("D5") "215.12"
$ws.Range("E5").Value = "  -1.42%  "

Set-TextValue $ws.Range("D6") "0.5054"
$ws.Range("E6").Value = "  -2.58%  "

$ws.Range("E7").Value = "  +0.11%  "

Set-TextValue $ws.Range("D8") "0.2574"
$ws.Range("E8").Value = "  +0.08%  "

Set-TextValue $ws.Range("D9") "0.06433"
$ws.Range("E9").Value = "  +0.35%  "

Set-TextValue $ws.Range("D10") "19.45"
$ws.Range("E10").Value = "  -2.28%  "

Set-TextValue $ws.Range("D11") "0.07730"

Set-TextValue $ws.Range("D12") "1.641.58"
$ws.Range("E12").Value = "  -1.07%  "

Set-TextValue $ws.Range("D13") "4.247"
$ws.Range("E13").Value = "  -1.23%  "

Set-TextValue $ws.Range("D14") "1.865.79"

Set-TextValue $ws.Range("D15") "0.5446"
$ws.Range("E15").Value = "  -1.49%  "

Set-TextValue $ws.Range("D16") "0.0₅7894"
$ws.Range("E16").Value = "  -1.89%  "

Set-TextValue $ws.Range("D17") "63.62"
$ws.Range("E17").Value = "  -1.12%  "

Set-TextValue $ws.Range("D18") "26.031.01"
$ws.Range("E18").Value = "  -0.75%  "

Set-TextValue $ws.Range("D19") "1.009"
$ws.Range("E19").Value = "  -0.07%  "

Set-TextValue $ws.Range("D20") "203.89"

Set-TextValue $ws.Range("D21") "4.290"
$ws.Range("E21").Value = "  -2.04%  "

Set-TextValue $ws.Range("D22") "9.982"
$ws.Range("E22").Value = "  -0.96%  "

Set-TextValue $ws.Range("D23") "5.952"

Set-TextValue $ws.Range("D24") "1.009"
$ws.Range("E24").Value = "  +0.13%  "

Set-TextValue $ws.Range("D25") "1.927"
$ws.Range("E25").Value = "  +9.36%  "

Set-TextValue $ws.Range("D26") "141.45"
$ws.Range("E26").Value = "  -1.68%  "

Set-TextValue $ws.Range("D27") "0.1154"
$ws.Range("E27").Value = "  -0.94%  "

Set-TextValue $ws.Range("D28") "15.74"
$ws.Range("E28").Value = "  -0.16%  "

Set-TextValue $ws.Range("D29") "6.729"
$ws.Range("E29").Value = "  -3.37%  "

Set-TextValue $ws.Range("D30") "0.05055"
$ws.Range("E30").Value = "  -4.25%  "

$ws.Range("E31").Value = "  -1.24%  "

Set-TextValue $ws.Range("D32") "3.248"
$ws.Range("E32").Value = "  -3.31%  "

Set-TextValue $ws.Range("D33") "3.192"
$ws.Range("E33").Value = "  -0.83%  "

Set-TextValue $ws.Range("D34") "1.542"
$ws.Range("E34").Value = "  -2.00%  "

Set-TextValue $ws.Range("D35") "2.339"
$ws.Range("E35").Value = "  -0.98%  "

Set-TextValue $ws.Range("D36") "2.629"
$ws.Range("E36").Value = "  -4.92%  "

Set-TextValue $ws.Range("D37") "0.8904"
$ws.Range("E37").Value = "  -3.91%  "

Set-TextValue $ws.Range("D38") "0.5616"
$ws.Range("E38").Value = "  -1.69%  "

Set-TextValue $ws.Range("D39") "1.145.92"
$ws.Range("E39").Value = "  -1.49%  "

$ws.Range("E40").Value = "  -1.40%  "

Set-TextValue $ws.Range("D41") "2.567"
$ws.Range("E41").Value = "  -0.53%  "

$ws.Range("E42").Value = "  +0.14%  "

Set-TextValue $ws.Range("D43") "5.657"
$ws.Range("E43").Value = "  -0.24%  "

Set-TextValue $ws.Range("D44") "0.8094"
$ws.Range("E44").Value = "  -3.67%  "

Set-TextValue $ws.Range("D45") "99.85"
$ws.Range("E45").Value = "  -0.03%  "

Set-TextValue $ws.Range("D46") "1.777.77"
$ws.Range("E46").Value = "  -1.08%  "

Set-TextValue $ws.Range("D47") "0.0₈113"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("E48").Value = "  +0.45%  "

$ws.Range("E49").Value = "  -0.33%  "

Set-TextValue $ws.Range("D50") "54.90"
$ws.Range("E50").Value = "  -2.06%  "

$ws.Range("E51").Value = "  -1.05%  "
